$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Duplicate the "SMALL MAGNET" block (rows 10-17) down to a new
#    "SPEHRICAL MAGNET" block (rows 19-26), reusing existing cell
#    formats/merges instead of minting brand-new ones.
# ------------------------------------------------------------------

# First pass: paste everything (values + formats + merged-cell shape)
# so the row/merge structure (including the two blank rows) is created.
$ws.Range("A10:G17").Copy()
$ws.Range("A19").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# Second pass: paste just the formats again so the cells pick back up
# the SAME style indexes the "SMALL MAGNET" block already uses (rather
# than the ad-hoc ones the first paste created).
$ws.Range("A10:G17").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row heights / thick-bottom banding to match rows 9/10/17 -> 18/19/26
$ws.Rows(18).RowHeight = 15.75
$ws.Rows(19).RowHeight = 15.75
$ws.Rows(26).RowHeight = 15.75

# ------------------------------------------------------------------
# 2. New header text for the duplicated block.
# ------------------------------------------------------------------
$ws.Range("A19").Value = "SPEHRICAL MAGNET"

# ------------------------------------------------------------------
# 3. New data values for the spherical-magnet measurements.
# ------------------------------------------------------------------
$ws.Range("A21").Value = 0.995
$ws.Range("B21").Value = 0.075
$ws.Range("E21").Formula = "=0.998^2"
$ws.Range("F21").Value = 0.075

$ws.Range("A22").Value = 0.64
$ws.Range("B22").Value = 0.085
$ws.Range("E22").Formula = "=0.796^2"
$ws.Range("F22").Value = 0.085

$ws.Range("A23").Value = 0.435
$ws.Range("B23").Value = 0.095
$ws.Range("E23").Formula = "=0.6525^2"
$ws.Range("F23").Value = 0.095

# Formulas (kept per-cell/relative; functionally identical to the
# shared formulas Excel would otherwise store).
$ws.Range("C21").Formula = "=A21*B21^6"
$ws.Range("C22").Formula = "=A22*B22^6"
$ws.Range("C23").Formula = "=A23*B23^6"

$ws.Range("G21").Formula = "=E21*F21^6"
$ws.Range("G22").Formula = "=E22*F22^6"
$ws.Range("G23").Formula = "=E23*F23^6"

$ws.Range("C26").Formula = "=AVERAGE(C21:C23)"
$ws.Range("G26").Formula = "=AVERAGE(G21:G23)"

# ------------------------------------------------------------------
# 4. Selection / view bookkeeping (mirrors the saved selection in the
#    target file: scrolled down to the new block, E24 selected).
# ------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 23
$ws.Range("E24").Select() | Out-Null

Write-Host "done"
